$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("out_vars")
$ws.Range("A111:J111").Copy()
$ws.Range("A112:J112").Select()
$ws.Paste()
